# Apply weekly work-report update:
#  - refresh report metadata (generated date, totals, billing period)
#  - add 3 new "Trans" line items (SVC-10-TP-AAL-RS) to the detail table
#  - renumber / re-total the TOTAL row accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Report summary / header cell updates
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:47 AM"
$ws.Range("C8").Value = 36879.36
$ws.Range("C9").Value = 143
$ws.Range("C10").Value = "07/14/2025 to 07/20/25"

# ---------------------------------------------------------------------------
# 2. Insert 3 new blank rows at the positions where new line items belong.
#    (Row numbers below refer to the "live" sheet at the moment of each
#    insert, i.e. after the previous inserts have already shifted rows.)
# ---------------------------------------------------------------------------
$ws.Rows.Item(154).Insert()   # new row before former row154 (Point 03 / DEC-20AL / Inst)
$ws.Rows.Item(159).Insert()   # new row before former row158 (Point 06 / DEG-40-PNA / Inst), now at 159
$ws.Rows.Item(160).Insert()   # new row right after the previous insert

# ---------------------------------------------------------------------------
# 3. Re-apply the correct alternating row style (the detail table alternates
#    a "white" style block (9/10/11) and a "shaded" style block (12/13/14)
#    strictly by row position) across the whole affected block, rows 154-163.
# ---------------------------------------------------------------------------
$templateWhite = $ws.Range("A16:H16")    # known-good style-group 9/10/11 row
$templateShade = $ws.Range("A17:H17")    # known-good style-group 12/13/14 row

for ($r = 154; $r -le 163; $r++) {
    if ((($r - 154) % 2) -eq 0) {
        $templateWhite.Copy()
    } else {
        $templateShade.Copy()
    }
    $ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Fill in the 3 new line-item rows
# ---------------------------------------------------------------------------
$ws.Range("A154").Value = "Point 03"
$ws.Range("B154").Value = "SVC-10-TP-AAL-RS"
$ws.Range("C154").Value = "Trans"
$ws.Range("D154").Value = "SVC,1/0,Trip,All Alum,Res"
$ws.Range("E154").Value = "EA"
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = ""
$ws.Range("H154").Value = 195.83

$ws.Range("A159").Value = "Point 04"
$ws.Range("B159").Value = "SVC-10-TP-AAL-RS"
$ws.Range("C159").Value = "Trans"
$ws.Range("D159").Value = "SVC,1/0,Trip,All Alum,Res"
$ws.Range("E159").Value = "EA"
$ws.Range("F159").Value = 1
$ws.Range("G159").Value = ""
$ws.Range("H159").Value = 195.83

$ws.Range("A160").Value = "Point 06"
$ws.Range("B160").Value = "SVC-10-TP-AAL-RS"
$ws.Range("C160").Value = "Trans"
$ws.Range("D160").Value = "SVC,1/0,Trip,All Alum,Res"
$ws.Range("E160").Value = "EA"
$ws.Range("F160").Value = 1
$ws.Range("G160").Value = ""
$ws.Range("H160").Value = 195.83

# ---------------------------------------------------------------------------
# 5. Update TOTAL row (now at row 164) with the new grand total
# ---------------------------------------------------------------------------
$ws.Range("H164").Value = 5805.479999999998

Write-Host "Edit complete"
